$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.518000000000001
$ws.Range("D6").Value = -7.388
$ws.Range("D7").Value = -7.258999999999999
$ws.Range("B8").Value = 6.359
$ws.Range("D8").Value = -7.419
$ws.Range("A12").Value = -21.401
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 6.679
$ws.Range("D19").Value = -8.059999999999999
$ws.Range("D21").Value = -7.25
$ws.Range("B22").Value = 6.461
$ws.Range("D24").Value = -7.431999999999999
